$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.32689766666667
$ws.Range("H2").Value = 30.980693
$ws.Range("I2").Value = 0.2044815006034941
$ws.Range("J2").Value = 0.204481500603494
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.309350333333333
$ws.Range("N2").Value = 24.928051
$ws.Range("O2").Value = 0.1535033474258946
$ws.Range("P2").Value = 0.1535033474258946
$ws.Range("Q2").Value = 85.80981056881588
$ws.Range("R2").Value = 772.2882951193429
$ws.Range("S2").Value = 0.03138859482930642
$ws.Range("T2").Value = 0.03138859482930642
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.32689766666667
$ws.Range("H3").Value = 30.980693
$ws.Range("I3").Value = 0.2044815006034941
$ws.Range("J3").Value = 0.204481500603494
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 37.153391
$ws.Range("N3").Value = 111.460173
$ws.Range("O3").Value = 0.6863556906301786
$ws.Range("P3").Value = 0.6863556906301786
$ws.Range("Q3").Value = 383.6792668266543
$ws.Range("R3").Value = 3453.113401439889
$ws.Range("S3").Value = 0.1403470415678065
$ws.Range("T3").Value = 0.1403470415678064
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.32689766666667
$ws.Range("H4").Value = 30.980693
$ws.Range("I4").Value = 0.2044815006034941
$ws.Range("J4").Value = 0.204481500603494
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.668653666666666
$ws.Range("N4").Value = 26.005961
$ws.Range("O4").Value = 0.1601409619439267
$ws.Range("P4").Value = 0.1601409619439267
$ws.Range("Q4").Value = 89.52029932344142
$ws.Range("R4").Value = 805.6826939109729
$ws.Range("S4").Value = 0.03274586420638118
$ws.Range("T4").Value = 0.03274586420638118
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.37031933333333
$ws.Range("H5").Value = 88.110958
$ws.Range("I5").Value = 0.5815577111671272
$ws.Range("J5").Value = 0.5815577111671272
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.309350333333333
$ws.Range("N5").Value = 24.928051
$ws.Range("O5").Value = 0.1535033474258946
$ws.Range("P5").Value = 0.1535033474258946
$ws.Range("Q5").Value = 244.0482727425398
$ws.Range("R5").Value = 2196.434454682858
$ws.Range("S5").Value = 0.08927105538549558
$ws.Range("T5").Value = 0.0892710553854956
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 29.37031933333333
$ws.Range("H6").Value = 88.110958
$ws.Range("I6").Value = 0.5815577111671272
$ws.Range("J6").Value = 0.5815577111671272
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 37.153391
$ws.Range("N6").Value = 111.460173
$ws.Range("O6").Value = 0.6863556906301786
$ws.Range("P6").Value = 0.6863556906301786
$ws.Range("Q6").Value = 1091.206957986192
$ws.Range("R6").Value = 9820.862621875733
$ws.Range("S6").Value = 0.3991554444894195
$ws.Range("T6").Value = 0.3991554444894195
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 29.37031933333333
$ws.Range("H7").Value = 88.110958
$ws.Range("I7").Value = 0.5815577111671272
$ws.Range("J7").Value = 0.5815577111671272
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.668653666666666
$ws.Range("N7").Value = 26.005961
$ws.Range("O7").Value = 0.1601409619439267
$ws.Range("P7").Value = 0.1601409619439267
$ws.Range("Q7").Value = 254.6011263800708
$ws.Range("R7").Value = 2291.410137420638
$ws.Range("S7").Value = 0.09313121129221205
$ws.Range("T7").Value = 0.09313121129221205
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.80562866666667
$ws.Range("H8").Value = 32.416886
$ws.Range("I8").Value = 0.2139607882293788
$ws.Range("J8").Value = 0.2139607882293788
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.309350333333333
$ws.Range("N8").Value = 24.928051
$ws.Range("O8").Value = 0.1535033474258946
$ws.Range("P8").Value = 0.1535033474258946
$ws.Range("Q8").Value = 89.78775416324287
$ws.Range("R8").Value = 808.0897874691859
$ws.Range("S8").Value = 0.03284369721109259
$ws.Range("T8").Value = 0.03284369721109259
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.80562866666667
$ws.Range("H9").Value = 32.416886
$ws.Range("I9").Value = 0.2139607882293788
$ws.Range("J9").Value = 0.2139607882293788
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 37.153391
$ws.Range("N9").Value = 111.460173
$ws.Range("O9").Value = 0.6863556906301786
$ws.Range("P9").Value = 0.6863556906301786
$ws.Range("Q9").Value = 401.4657468534753
$ws.Range("R9").Value = 3613.191721681278
$ws.Range("S9").Value = 0.1468532045729527
$ws.Range("T9").Value = 0.1468532045729527
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 10.80562866666667
$ws.Range("H10").Value = 32.416886
$ws.Range("I10").Value = 0.2139607882293788
$ws.Range("J10").Value = 0.2139607882293788
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.668653666666666
$ws.Range("N10").Value = 26.005961
$ws.Range("O10").Value = 0.1601409619439267
$ws.Range("P10").Value = 0.1601409619439267
$ws.Range("Q10").Value = 89.52029932344142
$ws.Range("R10").Value = 843.0322730574459
$ws.Range("S10").Value = 0.03426388644533352
$ws.Range("T10").Value = 0.03426388644533352